# ---------------------------------------------------------------------------
# Commit: "Sun, Apr 19, 2020 11:04:48 AM"
#
# Two independent changes land in this commit:
#
#   1. Slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") has a 2-column table
#      (the cash-flow glossary table). Its table style is switched from
#      {BFA5B35F-0780-4C22-9AC8-46F7AF7D826B} to
#      {D8D359EB-0141-4C0B-833F-F9BFA8E8B99D}.
#
#   2. The presentation's theme colour palette is switched from the
#      "Integral" scheme back to the stock "Office" scheme (dk2/lt2/accent*/
#      hlink/folHlink all change; dk1/lt1 stay black/white). Font scheme and
#      format scheme are untouched - only the 12 theme colours move.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- helper: build the packed BGR-in-int value PowerPoint's RGB() macro
#     produces from individual R/G/B byte components, so Colors(i).RGB=
#     round-trips to the same srgbClr hex we pass in.
function VbaRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# ---------------------------------------------------------------------------
# 1) Table style on slide 16, shape 3 (the graphicFrame holding the table)
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{D8D359EB-0141-4C0B-833F-F9BFA8E8B99D}")

# ---------------------------------------------------------------------------
# 2) Theme colours: Integral -> Office
# ---------------------------------------------------------------------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# index : scheme slot  : target "Office" colour
#   1   : dk1           000000  (unchanged)
#   2   : lt1           FFFFFF  (unchanged)
#   3   : dk2           44546A
#   4   : lt2           E7E6E6
#   5   : accent1       5B9BD5
#   6   : accent2       ED7D31
#   7   : accent3       A5A5A5
#   8   : accent4       FFC000
#   9   : accent5       4472C4
#  10   : accent6       70AD47
#  11   : hlink         0563C1
#  12   : folHlink      954F72
$colorScheme.Colors(1).RGB  = VbaRGB 0x00 0x00 0x00
$colorScheme.Colors(2).RGB  = VbaRGB 0xFF 0xFF 0xFF
$colorScheme.Colors(3).RGB  = VbaRGB 0x44 0x54 0x6A
$colorScheme.Colors(4).RGB  = VbaRGB 0xE7 0xE6 0xE6
$colorScheme.Colors(5).RGB  = VbaRGB 0x5B 0x9B 0xD5
$colorScheme.Colors(6).RGB  = VbaRGB 0xED 0x7D 0x31
$colorScheme.Colors(7).RGB  = VbaRGB 0xA5 0xA5 0xA5
$colorScheme.Colors(8).RGB  = VbaRGB 0xFF 0xC0 0x00
$colorScheme.Colors(9).RGB  = VbaRGB 0x44 0x72 0xC4
$colorScheme.Colors(10).RGB = VbaRGB 0x70 0xAD 0x47
$colorScheme.Colors(11).RGB = VbaRGB 0x05 0x63 0xC1
$colorScheme.Colors(12).RGB = VbaRGB 0x95 0x4F 0x72
